$d = $word.ActiveDocument

# Update the date/weekday heading
$d.Content.Find.Execute("2026-02-18 Wednesday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2026-02-19 Thursday", 2)

# Update the division problems in the table, addressed by row/column so that
# duplicate values (e.g. "70÷3=") are each mapped to their correct target.
$t = $d.Tables.Item(1)

$updates = @(
    @{Row=1;  Col=1; Text="60÷2="},
    @{Row=1;  Col=2; Text="91÷7="},
    @{Row=1;  Col=3; Text="99÷8="},
    @{Row=1;  Col=4; Text="12÷8="},
    @{Row=1;  Col=5; Text="95÷2="},

    @{Row=5;  Col=1; Text="17÷5="},
    @{Row=5;  Col=2; Text="56÷4="},
    @{Row=5;  Col=3; Text="84÷2="},
    @{Row=5;  Col=4; Text="58÷8="},
    @{Row=5;  Col=5; Text="79÷2="},

    @{Row=9;  Col=1; Text="65÷9="},
    @{Row=9;  Col=2; Text="41÷4="},
    @{Row=9;  Col=3; Text="17÷5="},
    @{Row=9;  Col=4; Text="40÷5="},
    @{Row=9;  Col=5; Text="51÷2="},

    @{Row=13; Col=1; Text="47÷7="},
    @{Row=13; Col=2; Text="85÷4="},
    @{Row=13; Col=3; Text="29÷2="},
    @{Row=13; Col=4; Text="80÷5="},
    @{Row=13; Col=5; Text="97÷7="},

    @{Row=17; Col=1; Text="14÷2="},
    @{Row=17; Col=2; Text="78÷7="},
    @{Row=17; Col=3; Text="35÷8="},
    @{Row=17; Col=4; Text="85÷7="},
    @{Row=17; Col=5; Text="15÷2="}
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $cell.Range.Text = $u.Text
}
